$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Copy row formatting ("borders/number-format" template) down onto the new
#    rows (263-271) before writing any values, mirroring the row-striping
#    that the underlying table style applies to appended Google-Forms rows.
#    Columns A:L are always present; M/N are copied per-cell only where the
#    new row actually needs them, so we don't leave stray blank <c> cells.
# ---------------------------------------------------------------------------
$xlPasteFormats = -4122

function CopyFormat($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy()
    $ws.Range($dstAddr).PasteSpecial($xlPasteFormats)
}

function CopyRowFormat($srcRow, $dstRow) {
    CopyFormat ("A" + $srcRow + ":L" + $srcRow) ("A" + $dstRow + ":L" + $dstRow)
}

CopyRowFormat 11 263
CopyRowFormat 12 264
CopyRowFormat 19 265
CopyRowFormat 12 266
CopyRowFormat 15 267
CopyRowFormat 12 268
CopyRowFormat 11 269
CopyRowFormat 24 270
CopyRowFormat 143 271

# M/N only where the new row actually has content there.
CopyFormat "M3"   "M263"
CopyFormat "N3"   "N265"
CopyFormat "M3"   "M269"
CopyFormat "M2"   "M270"
# Row 271 also needs the M (blank, bordered) and N (bordered, text) styles
# that only exist on the soon-to-be-cleared M262/N262 and on N191.
CopyFormat "M262" "M271"
CopyFormat "N191" "N271"
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Row 262 no longer is the last table row, so its trailing placeholder
#    cells (M262/N262, previously blank with a bottom border) go away.
# ---------------------------------------------------------------------------
$ws.Cells.Item(262, 13).Clear()
$ws.Cells.Item(262, 14).Clear()

# ---------------------------------------------------------------------------
# 3. Write the 9 new Google-Forms responses (rows 263-271).
# ---------------------------------------------------------------------------
$rows = @(
    @{ Row=263; A=45599.017438298615; B="lc990728@naver.com";      C="소프트웨어학부"; D=20245233; E="이하은"; F="대한민국"; G="취업자 / 15세 이상 인구"; H="조사대상 주간에 수입 있는 일을 하지 않았으나, 지난 4주간 구직활동하였으며 즉시 취업 가능한자"; I="평균 : 100만원, 중위값 : 200만원"; J="평균 : 100만원, 중위값 : 1,000만원"; K="""19.1%"""; L="Red";   M="나는 사후 장기기증에 참여할 뜻이 없다" },
    @{ Row=264; A=45599.030809004631; B="ggamy3637@naver.com";     C="러시아학과";     D=20201733; E="최효윤"; F="대한민국"; G="취업자 / 15세 이상 인구"; H="조사대상 주간에 수입 있는 일을 하지 않았으나, 지난 4주간 구직활동하였으며 즉시 취업 가능한자"; I="평균 : 200만원, 중위값 : 100만원"; J="평균 : 1,000만원, 중위값 : 100만원"; K="""19.1%"""; L="Black" },
    @{ Row=265; A=45599.045796226856; B="yumi0901gami@gmail.com";  C="중국학과";       D=20241515; E="김은주"; F="대한민국"; G="경제활동인구 / 15세이상 인구"; H="조사대상 주간에 일시적으로 병이 나거나, 날씨가 안 좋거나, 휴가 또는 연가, 노동쟁의 등의 이유로 일을 못한 일시 휴직자"; I="평균 : 100만원, 중위값 : 200만원"; J="평균 : 200만원, 중위값 : 100만원"; K="""15%"""; L="Black"; N="나는 사후 장기기증에 참여할 뜻이 있다" },
    @{ Row=266; A=45599.051838518513; B="snp040609@naver.com";     C="경영학과";       D=20242957; E="박세나"; F="스페인";   G="취업자 / 15세 이상 인구"; H="조사 대상 주간 중 수입을 목적으로 1시간 이상 일한 자"; I="평균 : 200만원, 중위값 : 100만원"; J="평균 : 1,000만원, 중위값 : 100만원"; K="""19.1%"""; L="Black" },
    @{ Row=267; A=45599.054907939819; B="hsjenny99@gmail.com";     C="소프트웨어학부"; D=20245246; E="전소현"; F="대한민국"; G="취업자 / 15세 이상 인구"; H="조사대상 주간에 수입 있는 일을 하지 않았으나, 지난 4주간 구직활동하였으며 즉시 취업 가능한자"; I="평균 : 200만원, 중위값 : 100만원"; J="평균 : 1,000만원, 중위값 : 100만원"; K="""19.1%"""; L="Red" },
    @{ Row=268; A=45599.056383969903; B="kbi70722@gmail.com";      C="일본학과";       D=20191604; E="김병일"; F="대한민국"; G="실업자 / 경제활동인구"; H="자기 가구에서 경영하는 농장이나 사업체의 수입을 높이는 데 도운 가족종사자로서 주당 18시간 이상 일한 자"; I="평균 : 100만원, 중위값 : 200만원"; J="평균 : 100만원, 중위값 : 1,000만원"; K="""15%"""; L="Red" },
    @{ Row=269; A=45599.102343344908; B="jina20050429@gmail.com";  C="환경생명공학과"; D=20243731; E="이진아"; F="OECD 전체"; G="경제활동인구 / 15세이상 인구"; H="조사대상 주간에 일시적으로 병이 나거나, 날씨가 안 좋거나, 휴가 또는 연가, 노동쟁의 등의 이유로 일을 못한 일시 휴직자"; I="평균 : 200만원, 중위값 : 100만원"; J="평균 : 1,000만원, 중위값 : 100만원"; K="""10%"""; L="Red"; M="나는 사후 장기기증에 참여할 뜻이 없다" },
    @{ Row=270; A=45599.118382175926; B="gangjunu@naver.com";      C="금융재무학과";   D=20242901; E="강준우"; F="대한민국"; G="취업자 / 경제활동인구"; H="조사 대상 주간 중 수입을 목적으로 1시간 이상 일한 자"; I="평균 : 200만원, 중위값 : 100만원"; J="평균 : 1,000만원, 중위값 : 100만원"; K="""19.1%"""; L="Red"; M="나는 사후 장기기증에 참여할 뜻이 없다" },
    @{ Row=271; A=45599.219025648148; B="jyj111212@naver.com";     C="인문학부";       D=20241083; E="장예지"; F="OECD 전체"; G="경제활동인구 / 15세이상 인구"; H="조사대상 주간에 수입 있는 일을 하지 않았으나, 지난 4주간 구직활동하였으며 즉시 취업 가능한자"; I="평균 : 200만원, 중위값 : 100만원"; J="평균 : 100만원, 중위값 : 1,000만원"; K="""5%""";  L="Black"; N="나는 사후 장기기증에 참여할 뜻이 있다" }
)

foreach ($r in $rows) {
    $rn = $r.Row
    $ws.Cells.Item($rn, 1).Value2 = $r.A
    $ws.Cells.Item($rn, 2).Value2 = $r.B
    $ws.Cells.Item($rn, 3).Value2 = $r.C
    $ws.Cells.Item($rn, 4).Value2 = $r.D
    $ws.Cells.Item($rn, 5).Value2 = $r.E
    $ws.Cells.Item($rn, 6).Value2 = $r.F
    $ws.Cells.Item($rn, 7).Value2 = $r.G
    $ws.Cells.Item($rn, 8).Value2 = $r.H
    $ws.Cells.Item($rn, 9).Value2 = $r.I
    $ws.Cells.Item($rn, 10).Value2 = $r.J
    $ws.Cells.Item($rn, 11).Value2 = $r.K
    $ws.Cells.Item($rn, 12).Value2 = $r.L
    if ($r.ContainsKey("M")) { $ws.Cells.Item($rn, 13).Value2 = $r.M }
    if ($r.ContainsKey("N")) { $ws.Cells.Item($rn, 14).Value2 = $r.N }
}

# ---------------------------------------------------------------------------
# 4. Grow the "Form_Responses1" table to cover the freshly written rows.
# ---------------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:N271"))

# ---------------------------------------------------------------------------
# 5. Move the UI cursor the same way the author's session ended up: scrolled
#    down a bit further and sitting on C281.
# ---------------------------------------------------------------------------
$ws.Range("C281").Select()
